# Applies the TEST-cyclic.xlsx change:
#  - adds a new worksheet "Sheet7" (OFFSET function reference sample) after "Sheet6"
#  - scrolls "Sheet6" so row 7 is the first visible row (view only)
#  - the newly added sheet becomes the active tab

$wb = $excel.ActiveWorkbook

# --- Sheet6: adjust the scroll position (best effort; selection itself is untouched) ---
$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# --- add the new worksheet as the last tab, after Sheet6 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $lastSheet)
$ws7.Name = "Sheet7"

# Row 5 - "Data" markers above the source table (filled in first, so "Data" is
# registered as a shared string before the header/description text below)
$ws7.Range("B5").Value = "Data"
$ws7.Range("C5").Value = "Data"

# Row 6-8 - source data table referenced by the OFFSET formulas below
$ws7.Range("B6").Value = 4
$ws7.Range("C6").Value = 10
$ws7.Range("B7").Value = 8
$ws7.Range("C7").Value = 3
$ws7.Range("B8").Value = 3
$ws7.Range("C8").Value = 6

# Header row
$ws7.Range("A1").Value = "Formula"
$ws7.Range("B1").Value = "Description"
$ws7.Range("C1").Value = "Result"

# Row 2 - OFFSET returning a single cell (B6)
$ws7.Range("A2").Formula = "=OFFSET(D3,3,-2,1,1)"
$ws7.Range("B2").Value = "Displays the value in cell B6 (4)"
$ws7.Range("C2").Value = 4

# Row 3 - SUM over an OFFSET range (B6:D8)
$ws7.Range("A3").Formula = "=SUM(OFFSET(D3:F5,3,-2, 3, 3))"
$ws7.Range("B3").Value = "Sums the range B6:D8"
$ws7.Range("C3").Value = 34

# Row 4 - OFFSET referencing a non-existent range -> #REF!
$ws7.Range("A4").Formula = "=OFFSET(D3, -3, -3)"
$ws7.Range("B4").Value = "Returns an error, because the reference is to a non-existent range on the worksheet."
$ws7.Range("C4").Value = "#REF!"

# Copy the "reference sample" cell style (font) from Sheet6!A7 onto A6, and match its row height
$ws6.Range("A7").Copy() | Out-Null
$ws7.Range("A6").PasteSpecial(-4122) | Out-Null
$ws7.Application.CutCopyMode = $false
$ws7.Rows.Item(6).RowHeight = 16.2

# Selection state that was captured for the new sheet (whole row 5 selected)
$ws7.Rows.Item(5).Select() | Out-Null
